# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column
# (value "stock" for the existing row), inserted between the existing
# "total" and "date" columns. Every column from "date" onward shifts
# one place to the right (date, legislator_name, legislator_id).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# "total" is column G (7); "date" is column H (8). Insert a fresh
# column at H so everything currently at/after H shifts right by one.
$ws.Columns.Item(8).Insert()

$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(2, 8).Value = "stock"
